$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the C1 header from "半价" to "班级" (D1 "性别" / D2 "女" remain as-is,
# their shared-string index merely shifts because "半价" is removed from the
# shared strings table).
$ws.Range("C1").Value = "班级"

# Select cell C1 to match the saved view state.
$ws.Range("C1").Select()
